# Update "想去人数" (interested-people count) figures for the gh-pages data refresh
# (commit: "Update gh-pages to output generated at 456a3b4")
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 894  # was 891
$ws.Range("F7").Value = 461  # was 460
$ws.Range("F9").Value = 2176  # was 2174
$ws.Range("F10").Value = 624  # was 623
$ws.Range("F13").Value = 1067  # was 1065
$ws.Range("F14").Value = 180  # was 178
$ws.Range("F15").Value = 2200  # was 2198
$ws.Range("F16").Value = 657  # was 654
$ws.Range("F17").Value = 12726  # was 12635
$ws.Range("F18").Value = 1238  # was 1237
$ws.Range("F19").Value = 12  # was 9
$ws.Range("F22").Value = 23  # was 22
$ws.Range("F25").Value = 265  # was 264

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4  # was 3
$ws.Range("F4").Value = 3  # was 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 478  # was 477
$ws.Range("F4").Value = 467  # was 466

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 478  # was 477
$ws.Range("F5").Value = 467  # was 466
$ws.Range("F7").Value = 4  # was 3
$ws.Range("F9").Value = 894  # was 891
$ws.Range("F10").Value = 3  # was 2
$ws.Range("F11").Value = 461  # was 460
$ws.Range("F13").Value = 2176  # was 2174
$ws.Range("F14").Value = 624  # was 623
$ws.Range("F19").Value = 1067  # was 1065
$ws.Range("F21").Value = 180  # was 178
$ws.Range("F24").Value = 2200  # was 2198
$ws.Range("F25").Value = 657  # was 654
$ws.Range("F28").Value = 1238  # was 1237
$ws.Range("F29").Value = 12  # was 9
$ws.Range("F32").Value = 23  # was 22
$ws.Range("F38").Value = 265  # was 264
